$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Numeric value updates ---
$ws.Range("M15").Value = -33.333333333333
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 20
$ws.Range("M16").Value = -9.43396226415
$ws.Range("N16").Value = -84
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 5
$ws.Range("H17").Value = -44.444444444444
$ws.Range("I17").Value = 67
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = 63.414634146341
$ws.Range("L17").Value = 86.111111111111
$ws.Range("M17").Value = 509.090909090909
$ws.Range("N17").Value = -22.988505747126
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = 76.086956521739
$ws.Range("L18").Value = 92.857142857142
$ws.Range("M18").Value = 15.714285714285
$ws.Range("N18").Value = -70.758122743682
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -2.702702702702
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 191
$ws.Range("K19").Value = 9.424083769633
$ws.Range("L19").Value = 36.601307189542
$ws.Range("M19").Value = 32.278481012658
$ws.Range("N19").Value = 26.666666666666
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 37.5
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = 57.142857142857
$ws.Range("L20").Value = 158.823529411765
$ws.Range("M20").Value = 51.724137931034
$ws.Range("N20").Value = -89.622641509434
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -19.047619047619
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = -6.849315068493
$ws.Range("I21").Value = 452
$ws.Range("J21").Value = 353
$ws.Range("K21").Value = 28.045325779036
$ws.Range("L21").Value = 59.717314487632
$ws.Range("M21").Value = 39.506172839506
$ws.Range("N21").Value = -64.325177584846
$ws.Range("M22").Value = -27.272727272727
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 20.833333333333
$ws.Range("F24").Value = 123
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 25.510204081632
$ws.Range("I24").Value = 617
$ws.Range("J24").Value = 454
$ws.Range("K24").Value = 35.90308370044
$ws.Range("L24").Value = 89.263803680981
$ws.Range("M24").Value = 102.960526315789
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 18
$ws.Range("H25").Value = -18.181818181818
$ws.Range("I25").Value = 113
$ws.Range("J25").Value = 102
$ws.Range("K25").Value = 10.78431372549
$ws.Range("L25").Value = 56.944444444444
$ws.Range("M25").Value = 48.684210526315
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -30
$ws.Range("N28").Value = -80
$ws.Range("N29").Value = -80

# --- Cells converting from numeric to suppressed-value text ("0" / "***.*") ---
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "0"
$ws.Range("G15").Copy()
$ws.Range("G26").PasteSpecial(-4122)

$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "***.*"
$ws.Range("H15").Copy()
$ws.Range("H26").PasteSpecial(-4122)

# --- Cells converting from suppressed-value text back to numeric ---
$ws.Range("D16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2

$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100

$excel.CutCopyMode = $false
